{"js": "// Pandoc-style title block: split \"On Pilgrimage - March/April 1975\" (Heading1)\n// and \"By Dorothy Day\" (bold run) into a Title paragraph and an Authors\n// paragraph, each made up of one run per token (word/punctuation/space),\n// matching the way pandoc emits docx title blocks.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst titleParagraph = paragraphs.items[0];\nconst authorParagraph = paragraphs.items[1];\n\n// Build a minimal OOXML package wrapping a single <w:p> whose runs are the\n// given tokens (each token becomes its own <w:r><w:t>).\nfunction buildParagraphOoxml(styleId, tokens) {\n  const runs = tokens\n    .map((token) => `<w:r><w:t xml:space=\"preserve\">${token}</w:t></w:r>`)\n    .join(\"\");\n  return (\n    '<?xml version=\"1.0\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    \"<pkg:xmlData>\" +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    \"<w:body>\" +\n    `<w:p><w:pPr><w:pStyle w:val=\"${styleId}\"/></w:pPr>${runs}</w:p>` +\n    \"</w:body>\" +\n    \"</w:document>\" +\n    \"</pkg:xmlData>\" +\n    \"</pkg:part>\" +\n    \"</pkg:package>\"\n  );\n}\n\nconst titleTokens = [\n  \"On\",\n  \" \",\n  \"Pilgrimage\",\n  \" \",\n  \"-\",\n  \" \",\n  \"March\",\n  \"/\",\n  \"April\",\n  \" \",\n  \"1975\",\n];\nconst authorTokens = [\"Dorothy\", \" \", \"Day\"];\n\n// Replace the author paragraph first so the (still-valid) titleParagraph\n// reference isn't affected by the sibling edit.\nauthorParagraph\n  .getRange()\n  .insertOoxml(buildParagraphOoxml(\"Authors\", authorTokens), Word.InsertLocation.replace);\nawait context.sync();\n\ntitleParagraph\n  .getRange()\n  .insertOoxml(buildParagraphOoxml(\"Title\", titleTokens), Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Pandoc-style title block: split \"On Pilgrimage - March/April 1975\" (Heading1)\n# and \"By Dorothy Day\" (bold run) into a Title paragraph and an Authors\n# paragraph, each made up of one run per token (word/punctuation/space),\n# matching the way pandoc emits docx title blocks.\n\n$d = $word.ActiveDocument\n\n$titleRange = $d.Paragraphs(1).Range\n$authorRange = $d.Paragraphs(2).Range\n\nfunction Build-ParagraphOoxml {\n    param(\n        [string]$StyleId,\n        [string[]]$Tokens\n    )\n    $runs = ($Tokens | ForEach-Object { '<w:r><w:t xml:space=\"preserve\">' + $_ + '</w:t></w:r>' }) -join \"\"\n    return '<?xml version=\"1.0\"?>' +\n        '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n        '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n        '<w:body>' +\n        '<w:p><w:pPr><w:pStyle w:val=\"' + $StyleId + '\"/></w:pPr>' + $runs + '</w:p>' +\n        '</w:body>' +\n        '</w:document>' +\n        '</pkg:xmlData>' +\n        '</pkg:part>' +\n        '</pkg:package>'\n}\n\n$titleTokens = @(\"On\", \" \", \"Pilgrimage\", \" \", \"-\", \" \", \"March\", \"/\", \"April\", \" \", \"1975\")\n$authorTokens = @(\"Dorothy\", \" \", \"Day\")\n\n# Replace the author paragraph's contents first (Range objects stay valid\n# across the sibling edit since no paragraphs are added or removed).\n$null = $authorRange.InsertXML((Build-ParagraphOoxml \"Authors\" $authorTokens))\n$null = $titleRange.InsertXML((Build-ParagraphOoxml \"Title\" $titleTokens))\n"}
